$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "66.179.67", "604.44").
# Force text format on the whole price column first so Excel COM does not
# auto-coerce these into numbers (which would drop formatting like trailing
# zeros, e.g. "1.00" -> 1) when we assign the new values below.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '66.179.67'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '3.602.17'
$ws.Range("E3").Value = '  +2.07%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '604.44'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").Value = '139.45'
$ws.Range("E6").Value = '  +1.46%  '
$ws.Range("D7").Value = '3.600.33'
$ws.Range("E7").Value = '  +2.10%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.500'
$ws.Range("E9").Value = '  +1.66%  '
$ws.Range("E10").Value = '  +2.76%  '
$ws.Range("E11").Value = '  +4.83%  '
$ws.Range("D12").Value = '0.394'
$ws.Range("E12").Value = '  +2.44%  '
$ws.Range("D13").Value = '4.220.15'
$ws.Range("E13").Value = '  +2.18%  '
$ws.Range("D14").Value = '28.45'
$ws.Range("E14").Value = '  +5.12%  '
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").Value = '3.607.34'
$ws.Range("E16").Value = '  +2.19%  '
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").Value = '66.280.06'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").Value = '10.17'
$ws.Range("E19").Value = '  +0.68%  '
$ws.Range("D20").Value = '14.69'
$ws.Range("E20").Value = '  +3.35%  '
$ws.Range("D21").Value = '5.91'
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("D22").Value = '397.46'
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("D23").Value = '0.591'
$ws.Range("E23").Value = '  +3.44%  '
$ws.Range("D24").Value = '3.750.88'
$ws.Range("E24").Value = '  +2.15%  '
$ws.Range("D25").Value = '75.09'
$ws.Range("E25").Value = '  +2.28%  '
$ws.Range("E27").Value = '  +6.59%  '
$ws.Range("E28").Value = '  +5.64%  '
$ws.Range("D29").Value = '1.66'
$ws.Range("E29").Value = '  +28.97%  '
$ws.Range("E30").Value = '  +7.12%  '
$ws.Range("E31").Value = '  +3.51%  '
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = '3.611.11'
$ws.Range("E33").Value = '  +1.65%  '
$ws.Range("D34").Value = '24.67'
$ws.Range("E34").Value = '  +3.80%  '
$ws.Range("E35").Value = '  +4.78%  '
$ws.Range("E37").Value = '  +9.09%  '
$ws.Range("E38").Value = '  +5.02%  '
$ws.Range("D39").Value = '7.07'
$ws.Range("E39").Value = '  +2.57%  '
$ws.Range("D40").Value = '168.49'
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("E41").Value = '  +5.81%  '
$ws.Range("D42").Value = '0.843'
$ws.Range("E42").Value = '  +2.46%  '
$ws.Range("E43").Value = '  +7.81%  '
$ws.Range("D44").Value = '26.19'
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").Value = '43.22'
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("E46").Value = '  +3.78%  '
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("E48").Value = '  +3.94%  '
$ws.Range("D49").Value = '7.02'
$ws.Range("E49").Value = '  +3.60%  '
$ws.Range("D50").Value = '2.460.68'
$ws.Range("E50").Value = '  +3.25%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").Value = '0.917'
$ws.Range("E51").Value = '  +10.70%  '

# Restore the default (no explicit style) formatting on column D now that
# the text values are safely stored, so styling matches the original file.
$priceRange.ClearFormats()
